$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move header labels (row 4) left into contiguous columns A:E
$ws.Range("B4").Cut($ws.Range("A4"))
$ws.Range("F4").Cut($ws.Range("B4"))
$ws.Range("J4").Cut($ws.Range("C4"))
$ws.Range("N4").Cut($ws.Range("D4"))
$ws.Range("R4").Cut($ws.Range("E4"))

# Move the data columns (rows 5:29) left into contiguous columns C:E
$ws.Range("J5:J29").Cut($ws.Range("C5"))
$ws.Range("N5:N29").Cut($ws.Range("D5"))
$ws.Range("R5:R29").Cut($ws.Range("E5"))

# Update the selected cell to match the new view
$ws.Range("I13").Select()
